$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the now-superfluous trailing empty row (the table shrinks
#     from A1:F14 to A1:F13) ---
$ws.Rows.Item(14).Delete()

# --- Row 11: new "Bowes et al., 2015" entry ---
$ws.Range("A11").Value = "Bowes et al., 2015"
$ws.Range("B11").Value = "British, Irish and Australians"
$ws.Range("C11").Value = "1,962/8,923"
$ws.Range("D11").Value = "PsA"
$ws.Range("E11").Value = "PsA specific 5q31 association"

# --- Row 12: new "Stuart et al., 2015" entry ---
$ws.Range("A12").Value = "Stuart et al., 2015"
$ws.Range("B12").Value = "White North American and European"
$ws.Range("C12").Value = "1,430/1,417"
$ws.Range("D12").Value = "PsA"
$ws.Range("E12").Value = " PsA versus psoriasis chr18 LOC100505817, psoriasis only RGS6"
$ws.Range("F12").Value = "Fine-mapping included and additional meta-analysis including psoriaisis"

# --- Row 13: previous "Tsoi et al., 2017*" row, now fully filled in ---
$ws.Range("A13").Value = "Tsoi et al., 2017*"
$ws.Range("B13").Value = "White North American and European"
$ws.Range("C13").Value = "19,032/39,498"
$ws.Range("D13").Value = "Psoriasis and PsA"
$ws.Range("E13").Value = "CHUK, IKBKE, FASLG,KLRK1,PTEN"
$ws.Range("F13").Value = "Largest meta-analysis so far"

# --- Apply the "light border, no fill" look used on the new highlighted
#     cells (matches the cellXfs entry that gained applyFill="1") ---
$ws.Range("A11:E11").Interior.ColorIndex = -4142
$ws.Range("A12").Interior.ColorIndex = -4142
$ws.Range("C12").Interior.ColorIndex = -4142
$ws.Range("E12").Interior.ColorIndex = -4142

# --- Column widths / row heights adjusted slightly as the table was
#     rebalanced ---
$ws.Columns.Item(1).ColumnWidth = 19.25
$ws.Columns.Item(2).ColumnWidth = 32.25
$ws.Columns.Item(3).ColumnWidth = 19.92
$ws.Columns.Item(4).ColumnWidth = 18.75
$ws.Columns.Item(6).ColumnWidth = 55.25

$ws.StandardHeight = 15

# --- Selection ends up on E12, matching the author's last edit ---
$ws.Range("E12").Select()
